$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.787.13'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '3.060.72'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.94'
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.96'
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.13%  '
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.376'
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("D12").Value = '3.587.38'
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("E13").Value = '  +2.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.33'
$ws.Range("E14").Value = '  +4.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000163'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").Value = '57.793.66'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").Value = '3.062.41'
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.11'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.82'
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.07'
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.85'
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.500'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.36'
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("E25").Value = '  +3.17%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '0.0₃0900'
$ws.Range("E27").Value = '  -4.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.45'
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.24'
$ws.Range("E29").Value = '  +6.70%  '
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("E31").Value = '  +2.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.64'
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.87'
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.96'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  +3.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0679'
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("D39").Value = '3.101.83'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.92'
$ws.Range("E40").Value = '  +3.03%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.655'
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").Value = '2.273.48'
$ws.Range("E44").Value = '  +3.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0258'
$ws.Range("E45").Value = '  +6.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.37'
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.70'
$ws.Range("E47").Value = '  +5.05%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.93'
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.936'
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.743'
$ws.Range("E50").Value = '  +9.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '255.18'
$ws.Range("E51").Value = '  +10.78%  '
